# Applies the "HITO: Funcionó prueba con Juan" edit:
#  - Row 2 (TKT-001) date moves 46016 -> 46017
#  - Row 3 (TKT-002) date moves 46016 -> 46017
#  - A brand-new row "TKT-M1-01" (Juan Perez) is inserted as the new row 4
#  - The old row 4 (TKT-003) becomes row 5, with its date bumped 46017 -> 46018
#  - A new audit row "TKT-OLD-01" (Juan Perez) is appended as row 6, dated 46016

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the dates for the first two existing rows.
$ws.Range("A2").Value = 46017
$ws.Range("A3").Value = 46017

# 2) Insert a brand new row above the current row 4 (TKT-003), pushing it down to row 5.
#    This keeps formatting (date style, etc.) consistent with the row being pushed down.
$ws.Rows.Item(4).Insert()

# 3) Fill in the new row 4 with the "TKT-M1-01" ticket for Juan Perez.
$ws.Range("A4").Value = 46018
$ws.Range("B4").Value = "TKT-M1-01"
$ws.Range("C4").Value = "Juan Perez"
$ws.Range("D4").Value = "XY-9999"
$ws.Range("E4").Value = "Cliente Manana"
$ws.Range("F4").Value = "Ruta 66 km 10"
$ws.Range("G4").Value = "Mantenimiento"

# 4) The old TKT-003 row, now row 5: reassign tech to Juan Perez and bump its date.
$ws.Range("A5").Value = 46018
$ws.Range("C5").Value = "Juan Perez"

# 5) Append a new audit row (row 6) at the bottom of the table.
#    Give A6 the same date number-format as the rows above it (style index 2).
$ws.Range("A6").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("A6").Value = 46016
$ws.Range("B6").Value = "TKT-OLD-01"
$ws.Range("C6").Value = "Juan Perez"
$ws.Range("D6").Value = "ZZ-0000"
$ws.Range("E6").Value = "Old Task"
$ws.Range("F6").Value = "Somewhere"
$ws.Range("G6").Value = "Audit"
